# Add season-record columns (Wins / Losses / Ties) to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: copy the existing header formatting (bold, centered, bordered)
# from A1 onto the three new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Data rows 2-45: every player on the roster shares the team's season record.
for ($row = 2; $row -le 45; $row++) {
    $ws.Cells.Item($row, 30).Value = 86
    $ws.Cells.Item($row, 31).Value = 76
    $ws.Cells.Item($row, 32).Value = 0
}
